# Add recording links to schedule
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain-text recording links for rows that only display link text
# (same visual treatment as the pre-existing blank "Hyperlink"-styled
# placeholder cells, or no special style where none existed before).
# Order matches row order top-to-bottom so shared-string indices line up
# with the source edit.
$ws.Range("K3").Value  = "https://youtu.be/OccS8ZFN0E4"

# K5 becomes a real, clickable hyperlink (new relationship + hyperlink entry).
$ws.Hyperlinks.Add($ws.Range("K5"), "https://youtu.be/gAAs4BWfAtM", "", "", "https://youtu.be/gAAs4BWfAtM")

$ws.Range("K10").Value = "https://youtu.be/23C5LFXZZ2U"
$ws.Range("K12").Value = "https://youtu.be/X6zQ0a0da8I"
$ws.Range("K16").Value = "https://youtu.be/a-vzYSYA9mg"
$ws.Range("K19").Value = "https://youtu.be/7D6iMgvPFZk"
$ws.Range("K24").Value = "https://youtu.be/FvIsISt8Auw"
$ws.Range("K25").Value = "https://youtu.be/-DnF0uD3euo"
$ws.Range("K27").Value = "https://youtu.be/lGAcQqeEOzA"
$ws.Range("K32").Value = "https://youtu.be/h8XZzdBMokc"
$ws.Range("K36").Value = "https://youtu.be/epN5KSt6qZg"

# Leave the cursor where the author ended up after adding the links.
$ws.Range("K42").Select()
